$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 1.893344666666666
$ws.Cells.Item(2, 8).Value2 = 5.680033999999999
$ws.Cells.Item(2, 9).Value2 = 0.05525983881677096
$ws.Cells.Item(2, 10).Value2 = 0.05525983881677096
$ws.Cells.Item(2, 13).Value2 = 1.815761
$ws.Cells.Item(2, 14).Value2 = 5.447283000000001
$ws.Cells.Item(2, 15).Value2 = 0.07007596730428067
$ws.Cells.Item(2, 16).Value2 = 0.07007596730428067
$ws.Cells.Item(2, 17).Value2 = 3.437861405291333
$ws.Cells.Item(2, 18).Value2 = 30.940752647622
$ws.Cells.Item(2, 19).Value2 = 0.003872386658163861
$ws.Cells.Item(2, 20).Value2 = 0.003872386658163861

# Row 3
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 1.893344666666666
$ws.Cells.Item(3, 8).Value2 = 5.680033999999999
$ws.Cells.Item(3, 9).Value2 = 0.05525983881677096
$ws.Cells.Item(3, 10).Value2 = 0.05525983881677096
$ws.Cells.Item(3, 15).Value2 = 0.5079540516959071
$ws.Cells.Item(3, 16).Value2 = 0.5079540516959072
$ws.Cells.Item(3, 17).Value2 = 24.91975062440622
$ws.Cells.Item(3, 18).Value2 = 224.277755619656
$ws.Cells.Item(3, 19).Value2 = 0.02806945902304157
$ws.Cells.Item(3, 20).Value2 = 0.02806945902304157

# Row 4
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 1.893344666666666
$ws.Cells.Item(4, 8).Value2 = 5.680033999999999
$ws.Cells.Item(4, 9).Value2 = 0.05525983881677096
$ws.Cells.Item(4, 10).Value2 = 0.05525983881677096
$ws.Cells.Item(4, 13).Value2 = 9.711409333333334
$ws.Cells.Item(4, 14).Value2 = 29.134228
$ws.Cells.Item(4, 15).Value2 = 0.3747940411327002
$ws.Cells.Item(4, 16).Value2 = 0.3747940411327002
$ws.Cells.Item(4, 17).Value2 = 18.38704506708356
$ws.Cells.Item(4, 18).Value2 = 165.483405603752
$ws.Cells.Item(4, 19).Value2 = 0.02071105830247923
$ws.Cells.Item(4, 20).Value2 = 0.02071105830247923

# Row 5
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 1.893344666666666
$ws.Cells.Item(5, 8).Value2 = 5.680033999999999
$ws.Cells.Item(5, 9).Value2 = 0.05525983881677096
$ws.Cells.Item(5, 10).Value2 = 0.05525983881677096
$ws.Cells.Item(5, 13).Value2 = 1.222391
$ws.Cells.Item(5, 14).Value2 = 3.667173
$ws.Cells.Item(5, 15).Value2 = 0.04717593986711188
$ws.Cells.Item(5, 16).Value2 = 0.04717593986711189
$ws.Cells.Item(5, 17).Value2 = 2.314407480431333
$ws.Cells.Item(5, 18).Value2 = 20.829667323882
$ws.Cells.Item(5, 19).Value2 = 0.002606934833086282
$ws.Cells.Item(5, 20).Value2 = 0.002606934833086282

# Row 6
$ws.Cells.Item(6, 7).Value2 = 4.159773333333334
$ws.Cells.Item(6, 9).Value2 = 0.1214086415227279
$ws.Cells.Item(6, 10).Value2 = 0.1214086415227279
$ws.Cells.Item(6, 13).Value2 = 1.815761
$ws.Cells.Item(6, 14).Value2 = 5.447283000000001
$ws.Cells.Item(6, 15).Value2 = 0.07007596730428067
$ws.Cells.Item(6, 16).Value2 = 0.07007596730428067
$ws.Cells.Item(6, 17).Value2 = 7.553154187506668
$ws.Cells.Item(6, 18).Value2 = 67.97838768756002
$ws.Cells.Item(6, 19).Value2 = 0.008507827993803813
$ws.Cells.Item(6, 20).Value2 = 0.008507827993803813

# Row 7
$ws.Cells.Item(7, 7).Value2 = 4.159773333333334
$ws.Cells.Item(7, 9).Value2 = 0.1214086415227279
$ws.Cells.Item(7, 10).Value2 = 0.1214086415227279
$ws.Cells.Item(7, 15).Value2 = 0.5079540516959071
$ws.Cells.Item(7, 16).Value2 = 0.5079540516959072
$ws.Cells.Item(7, 17).Value2 = 54.74994381409778
$ws.Cells.Item(7, 19).Value2 = 0.06167001137236559
$ws.Cells.Item(7, 20).Value2 = 0.0616700113723656

# Row 8
$ws.Cells.Item(8, 7).Value2 = 4.159773333333334
$ws.Cells.Item(8, 9).Value2 = 0.1214086415227279
$ws.Cells.Item(8, 10).Value2 = 0.1214086415227279
$ws.Cells.Item(8, 13).Value2 = 9.711409333333334
$ws.Cells.Item(8, 14).Value2 = 29.134228
$ws.Cells.Item(8, 15).Value2 = 0.3747940411327002
$ws.Cells.Item(8, 16).Value2 = 0.3747940411327002
$ws.Cells.Item(8, 17).Value2 = 40.39726157388445
$ws.Cells.Item(8, 18).Value2 = 363.5753541649601
$ws.Cells.Item(8, 19).Value2 = 0.04550323538473453
$ws.Cells.Item(8, 20).Value2 = 0.04550323538473453

# Row 9
$ws.Cells.Item(9, 7).Value2 = 4.159773333333334
$ws.Cells.Item(9, 9).Value2 = 0.1214086415227279
$ws.Cells.Item(9, 10).Value2 = 0.1214086415227279
$ws.Cells.Item(9, 13).Value2 = 1.222391
$ws.Cells.Item(9, 14).Value2 = 3.667173
$ws.Cells.Item(9, 15).Value2 = 0.04717593986711188
$ws.Cells.Item(9, 16).Value2 = 0.04717593986711189
$ws.Cells.Item(9, 17).Value2 = 5.084869484706667
$ws.Cells.Item(9, 18).Value2 = 45.76382536236
$ws.Cells.Item(9, 19).Value2 = 0.005727566771823954
$ws.Cells.Item(9, 20).Value2 = 0.005727566771823955

# Row 10
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 1.240212333333333
$ws.Cells.Item(10, 8).Value2 = 3.720637
$ws.Cells.Item(10, 9).Value2 = 0.03619728348733726
$ws.Cells.Item(10, 10).Value2 = 0.03619728348733727
$ws.Cells.Item(10, 13).Value2 = 1.815761
$ws.Cells.Item(10, 14).Value2 = 5.447283000000001
$ws.Cells.Item(10, 15).Value2 = 0.07007596730428067
$ws.Cells.Item(10, 16).Value2 = 0.07007596730428067
$ws.Cells.Item(10, 17).Value2 = 2.251929186585667
$ws.Cells.Item(10, 18).Value2 = 20.267362679271
$ws.Cells.Item(10, 19).Value2 = 0.002536559654162425
$ws.Cells.Item(10, 20).Value2 = 0.002536559654162425

# Row 11
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 1.240212333333333
$ws.Cells.Item(11, 8).Value2 = 3.720637
$ws.Cells.Item(11, 9).Value2 = 0.03619728348733726
$ws.Cells.Item(11, 10).Value2 = 0.03619728348733727
$ws.Cells.Item(11, 15).Value2 = 0.5079540516959071
$ws.Cells.Item(11, 16).Value2 = 0.5079540516959072
$ws.Cells.Item(11, 17).Value2 = 16.32337873398978
$ws.Cells.Item(11, 18).Value2 = 146.910408605908
$ws.Cells.Item(11, 19).Value2 = 0.01838655680777832
$ws.Cells.Item(11, 20).Value2 = 0.01838655680777833

# Row 12
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 1.240212333333333
$ws.Cells.Item(12, 8).Value2 = 3.720637
$ws.Cells.Item(12, 9).Value2 = 0.03619728348733726
$ws.Cells.Item(12, 10).Value2 = 0.03619728348733727
$ws.Cells.Item(12, 13).Value2 = 9.711409333333334
$ws.Cells.Item(12, 14).Value2 = 29.134228
$ws.Cells.Item(12, 15).Value2 = 0.3747940411327002
$ws.Cells.Item(12, 16).Value2 = 0.3747940411327002
$ws.Cells.Item(12, 17).Value2 = 12.04420962924844
$ws.Cells.Item(12, 18).Value2 = 108.397886663236
$ws.Cells.Item(12, 19).Value2 = 0.01356652615624509
$ws.Cells.Item(12, 20).Value2 = 0.01356652615624509

# Row 13
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 1.240212333333333
$ws.Cells.Item(13, 8).Value2 = 3.720637
$ws.Cells.Item(13, 9).Value2 = 0.03619728348733726
$ws.Cells.Item(13, 10).Value2 = 0.03619728348733727
$ws.Cells.Item(13, 13).Value2 = 1.222391
$ws.Cells.Item(13, 14).Value2 = 3.667173
$ws.Cells.Item(13, 15).Value2 = 0.04717593986711188
$ws.Cells.Item(13, 16).Value2 = 0.04717593986711189
$ws.Cells.Item(13, 17).Value2 = 1.516024394355667
$ws.Cells.Item(13, 18).Value2 = 13.644219549201
$ws.Cells.Item(13, 19).Value2 = 0.001707640869151425
$ws.Cells.Item(13, 20).Value2 = 0.001707640869151425

# Row 14
$ws.Cells.Item(14, 7).Value2 = 26.96925
$ws.Cells.Item(14, 8).Value2 = 80.90774999999999
$ws.Cells.Item(14, 9).Value2 = 0.7871342361731639
$ws.Cells.Item(14, 10).Value2 = 0.7871342361731638
$ws.Cells.Item(14, 13).Value2 = 1.815761
$ws.Cells.Item(14, 14).Value2 = 5.447283000000001
$ws.Cells.Item(14, 15).Value2 = 0.07007596730428067
$ws.Cells.Item(14, 16).Value2 = 0.07007596730428067
$ws.Cells.Item(14, 17).Value2 = 48.96971234925
$ws.Cells.Item(14, 18).Value2 = 440.72741114325
$ws.Cells.Item(14, 19).Value2 = 0.05515919299815058
$ws.Cells.Item(14, 20).Value2 = 0.05515919299815056

# Row 15
$ws.Cells.Item(15, 7).Value2 = 26.96925
$ws.Cells.Item(15, 8).Value2 = 80.90774999999999
$ws.Cells.Item(15, 9).Value2 = 0.7871342361731639
$ws.Cells.Item(15, 10).Value2 = 0.7871342361731638
$ws.Cells.Item(15, 15).Value2 = 0.5079540516959071
$ws.Cells.Item(15, 16).Value2 = 0.5079540516959072
$ws.Cells.Item(15, 17).Value2 = 354.962831839
$ws.Cells.Item(15, 18).Value2 = 3194.665486551
$ws.Cells.Item(15, 19).Value2 = 0.3998280244927216
$ws.Cells.Item(15, 20).Value2 = 0.3998280244927217

# Row 16
$ws.Cells.Item(16, 7).Value2 = 26.96925
$ws.Cells.Item(16, 8).Value2 = 80.90774999999999
$ws.Cells.Item(16, 9).Value2 = 0.7871342361731639
$ws.Cells.Item(16, 10).Value2 = 0.7871342361731638
$ws.Cells.Item(16, 13).Value2 = 9.711409333333334
$ws.Cells.Item(16, 14).Value2 = 29.134228
$ws.Cells.Item(16, 15).Value2 = 0.3747940411327002
$ws.Cells.Item(16, 16).Value2 = 0.3747940411327002
$ws.Cells.Item(16, 17).Value2 = 261.909426163
$ws.Cells.Item(16, 18).Value2 = 2357.184835467
$ws.Cells.Item(16, 19).Value2 = 0.2950132212892413
$ws.Cells.Item(16, 20).Value2 = 0.2950132212892413

# Row 17
$ws.Cells.Item(17, 7).Value2 = 26.96925
$ws.Cells.Item(17, 8).Value2 = 80.90774999999999
$ws.Cells.Item(17, 9).Value2 = 0.7871342361731639
$ws.Cells.Item(17, 10).Value2 = 0.7871342361731638
$ws.Cells.Item(17, 13).Value2 = 1.222391
$ws.Cells.Item(17, 14).Value2 = 3.667173
$ws.Cells.Item(17, 15).Value2 = 0.04717593986711188
$ws.Cells.Item(17, 16).Value2 = 0.04717593986711189
$ws.Cells.Item(17, 17).Value2 = 32.96696847675
$ws.Cells.Item(17, 18).Value2 = 296.7027162907499
$ws.Cells.Item(17, 19).Value2 = 0.03713379739305022
$ws.Cells.Item(17, 20).Value2 = 0.03713379739305023
